# Auto-generated Excel COM-interop script
# Applies the 2022-06-03 violent crime data update across all affected sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 2725
$ws.Range('I3').Value = 2851
$ws.Range('E4').Value = 1964
$ws.Range('F4').Value = 1853
$ws.Range('I4').Value = 695
$ws.Range('I6').Value = 3243
$ws.Range('E7').Value = 25968
$ws.Range('F7').Value = 24042
$ws.Range('I7').Value = 9765

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I7').Value = 324
$ws.Range('I8').Value = 622
$ws.Range('I10').Value = 72
$ws.Range('I11').Value = 161
$ws.Range('I20').Value = 248
$ws.Range('I21').Value = 58
$ws.Range('I24').Value = 24
$ws.Range('I29').Value = 652
$ws.Range('I31').Value = 90
$ws.Range('I33').Value = 456
$ws.Range('I34').Value = 42
$ws.Range('I36').Value = 130
$ws.Range('I37').Value = 316
$ws.Range('I42').Value = 338
$ws.Range('I44').Value = 74
$ws.Range('I48').Value = 108
$ws.Range('I49').Value = 69
$ws.Range('I52').Value = 204
$ws.Range('I53').Value = 113
$ws.Range('I54').Value = 221
$ws.Range('E63').Value = 311
$ws.Range('F63').Value = 150
$ws.Range('I63').Value = 35
$ws.Range('I64').Value = 90
$ws.Range('I65').Value = 216
$ws.Range('I66').Value = 23
$ws.Range('I67').Value = 377
$ws.Range('I71').Value = 23
$ws.Range('I73').Value = 83
$ws.Range('I76').Value = 152
$ws.Range('I78').Value = 130
$ws.Range('I79').Value = 247
$ws.Range('I82').Value = 9
$ws.Range('I83').Value = 196
$ws.Range('I84').Value = 78
$ws.Range('I85').Value = 451
$ws.Range('I88').Value = 82
$ws.Range('I90').Value = 113
$ws.Range('I91').Value = 116
$ws.Range('I94').Value = 89
$ws.Range('E101').Value = 25968
$ws.Range('F101').Value = 24042
$ws.Range('I101').Value = 9765

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I2').Value = 113
$ws.Range('I3').Value = 186
$ws.Range('I7').Value = 451

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I3').Value = 80
$ws.Range('I4').Value = 24
$ws.Range('I6').Value = 45
$ws.Range('I7').Value = 204

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('I2').Value = 72
$ws.Range('I7').Value = 161

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I3').Value = 173
$ws.Range('I6').Value = 196
$ws.Range('I7').Value = 622

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('I2').Value = 22
$ws.Range('I4').Value = 9
$ws.Range('I7').Value = 113

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('I6').Value = 83
$ws.Range('I7').Value = 324

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I2').Value = 102
$ws.Range('I6').Value = 87
$ws.Range('I7').Value = 316

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I2').Value = 86
$ws.Range('I3').Value = 132
$ws.Range('I6').Value = 133
$ws.Range('I7').Value = 377

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('I3').Value = 28
$ws.Range('I7').Value = 90

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('I6').Value = 19
$ws.Range('I7').Value = 78

$ws = $wb.Worksheets.Item('New City')
$ws.Range('I2').Value = 68
$ws.Range('I6').Value = 70
$ws.Range('I7').Value = 216

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I2').Value = 71
$ws.Range('I4').Value = 9
$ws.Range('I7').Value = 196

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I2').Value = 105
$ws.Range('I7').Value = 456

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('I6').Value = 40
$ws.Range('I7').Value = 69

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I3').Value = 47
$ws.Range('I6').Value = 107
$ws.Range('I7').Value = 221

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 203
$ws.Range('I6').Value = 178
$ws.Range('I7').Value = 652

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('I6').Value = 21
$ws.Range('I7').Value = 74

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I2').Value = 14
$ws.Range('I6').Value = 62
$ws.Range('I7').Value = 108

$ws = $wb.Worksheets.Item('River North')
$ws.Range('I6').Value = 61
$ws.Range('I7').Value = 152

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I2').Value = 89
$ws.Range('I3').Value = 113
$ws.Range('I6').Value = 92
$ws.Range('I7').Value = 338

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('I6').Value = 35
$ws.Range('I7').Value = 72

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I6').Value = 53
$ws.Range('I7').Value = 130

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('I2').Value = 7
$ws.Range('I7').Value = 24

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('I3').Value = 39
$ws.Range('I7').Value = 116

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('I2').Value = 4
$ws.Range('I7').Value = 58

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I2').Value = 70
$ws.Range('I7').Value = 247

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('I6').Value = 32
$ws.Range('I7').Value = 90

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I2').Value = 68
$ws.Range('I6').Value = 89
$ws.Range('I7').Value = 248

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I3').Value = 41
$ws.Range('I7').Value = 130

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('I2').Value = 16
$ws.Range('I7').Value = 42

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I3').Value = 15
$ws.Range('I7').Value = 89

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('I6').Value = 11
$ws.Range('I7').Value = 23

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I6').Value = 24
$ws.Range('I7').Value = 83

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I6').Value = 32
$ws.Range('I7').Value = 82

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('I3').Value = 21
$ws.Range('I7').Value = 113

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('I6').Value = 8
$ws.Range('I7').Value = 23

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range('I3').Value = 2
$ws.Range('I6').Value = 9
